$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.967.44"
$ws.Range("E2").Value = "  -1.19%  "

$ws.Range("D3").Value = "2.602.14"
$ws.Range("E3").Value = "  -1.80%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.15"
$ws.Range("E5").Value = "  +3.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.97"
$ws.Range("E6").Value = "  -1.66%  "

$ws.Range("E8").Value = "  +3.83%  "

$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("E10").Value = "  -1.52%  "

$ws.Range("E11").Value = "  +5.12%  "

$ws.Range("E12").Value = "  -0.84%  "

$ws.Range("D13").Value = "3.058.12"
$ws.Range("E13").Value = "  -2.04%  "

$ws.Range("D14").Value = "58.901.23"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.95"
$ws.Range("E15").Value = "  -2.24%  "

$ws.Range("D16").Value = "2.608.50"
$ws.Range("E16").Value = "  -2.58%  "

$ws.Range("E17").Value = "  -2.03%  "

$ws.Range("E18").Value = "  +0.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.29"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.13"
$ws.Range("E20").Value = "  -2.16%  "

$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.73"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("E24").Value = "  +2.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -1.97%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.19"
$ws.Range("E27").Value = "  -1.43%  "

$ws.Range("D28").Value = "0.0₃0757"
$ws.Range("E28").Value = "  +1.09%  "

$ws.Range("E30").Value = "  +1.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.02"
$ws.Range("E31").Value = "  +2.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.14"
$ws.Range("E32").Value = "  +2.07%  "

$ws.Range("E33").Value = "  +0.60%  "

$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.899"
$ws.Range("E35").Value = "  +6.79%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.883"
$ws.Range("E36").Value = "  +5.18%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.13"
$ws.Range("E37").Value = "  -0.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.94"
$ws.Range("E38").Value = "  -1.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "283.26"
$ws.Range("E41").Value = "  -1.07%  "

$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.602"
$ws.Range("E43").Value = "  -0.79%  "

$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0537"
$ws.Range("E45").Value = "  -0.53%  "

$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0229"
$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("D48").Value = "1.944.65"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.68"
$ws.Range("E49").Value = "  +5.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.13"
$ws.Range("E50").Value = "  -1.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.48"
$ws.Range("E51").Value = "  -1.92%  "
